# Roster update for Washington 2023 sheet:
#  - Row 17 (formerly "Jamaree Bouyea") is replaced by "Jay Huff (TW)"'s
#    data, and gains a jersey number (34) in column B.
#  - Row 18 (formerly "Jay Huff (TW)") is replaced by a new player,
#    "Xavier Cooks" (no jersey number, like the row he replaces).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 17: Jay Huff (TW) ---
$ws.Range("B17").Value = 34
$ws.Range("C17").Value = "Jay Huff (TW)"
$ws.Range("D17").Value = "C"
$ws.Range("E17").Value = "7-1"
$ws.Range("F17").Value = 240
$ws.Range("G17").Value = "August 25, 1998"
$ws.Range("H17").Value = "us"
# "Exp" of 1 needs to stay a text value (like "R" elsewhere in this column),
# so format the cell as Text first to stop COM's automatic number coercion.
$ws.Range("I17").NumberFormat = "@"
$ws.Range("I17").Value = "1"
$ws.Range("J17").Value = "Virginia"
$ws.Range("K17").Value = "https://www.basketball-reference.com/players/h/huffja01.html"

# --- Row 18: Xavier Cooks ---
$ws.Range("C18").Value = "Xavier Cooks"
$ws.Range("D18").Value = "PF"
$ws.Range("E18").Value = "6-8"
$ws.Range("F18").Value = 183
$ws.Range("G18").Value = "August 19, 1995"
$ws.Range("H18").Value = "au"
$ws.Range("I18").Value = "R"
$ws.Range("J18").Value = "Winthrop University"
$ws.Range("K18").Value = "https://www.basketball-reference.com/players/c/cooksxa01.html"
